# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated data (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F, applied identically to both sheets.
$updates = @{
    5  = 28
    7  = 2682
    9  = 1704
    12 = 573
    15 = 81
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
